# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy the existing header style (bold + border, style index
# used by A1:AC1) onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows (2-53): every player on the roster shares the team's season
# record, which was 69 wins, 93 losses, 0 ties.
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 69
    $ws.Cells.Item($r, 31).Value = 93
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-53"
